$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as a literal (shared-string) cell without ever
# letting Excel's "General format" auto-number-detection kick in and
# without touching the cell's existing style. We do this by putting a
# quoted-text FORMULA into the cell (so the result is always Text, even
# for digit-only codes), then collapsing the formula down to its
# computed value via Copy + PasteSpecial(values only). A plain
# Range.Value assignment would re-run the same "does this look like a
# number?" inference that created the need for this workaround, and
# forcing NumberFormat="@" (Text) would allocate a brand-new cell style.
function Set-TextValue($cell, $text) {
  $escaped = $text.Replace('"', '""')
  $ws.Range($cell).Formula = '="' + $escaped + '"'
  $ws.Range($cell).Copy()
  $ws.Range($cell).PasteSpecial(-4163)
}

# Remove the trailing rows (10-14) that no longer exist in the new table.
$ws.Rows("10:14").Delete()

# Column F narrows from 12 to 11 characters wide. ColumnWidth (COM units)
# runs 0.8333.. below the raw OOXML character width on this sheet's font,
# so ask for 11 - 5/6 to land exactly on width="11" after save.
$ws.Columns("F:F").ColumnWidth = 11 - (5/6)

# Replace the product rows (2-9) with the new data set.
$data = @(
  @("10002350","ABC KCP MANIS TGG275","TH1MKT","2","1","RT,(E-2B)"),
  @("20029222","PRONAS CORNED BEEF50","TH1MKT","2","2","RT,(E-4B)"),
  @("10013205","SUNLIGHT J/NIPIS.750","TH1MKT","3","1","PT"),
  @("20031238","WHISKAS JR MCKRL 80","TH1MKT","3","2","RT,(E-3B)"),
  @("20031233","WHISKAS CF MAC&SAL80","TH1MKT","3","3","RT,(E-3B)"),
  @("10037405","C/LANG KAYU PUTIH 30","TH1MKT","3","4","RT,(E-6B)"),
  @("20134103","ATTACK GEL S.CNTA515","TH1MKT","4","1","RT,(E-1B)"),
  @("20032250","KISPRAY VIOLET PC280","TH1MKT","4","2","RT")
)

$cols = @("A","B","C","D","E","F")
$r = 2
foreach ($row in $data) {
  for ($i = 0; $i -lt 6; $i++) {
    Set-TextValue ($cols[$i] + $r) $row[$i]
  }
  $r = $r + 1
}
